$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition listing) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3
$ws1.Range("F3").Value = 7

# Sheet "全部类型" (All types) contains the same data and needs the same update
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3
$ws4.Range("F3").Value = 7
